# "Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)"
# Refresh the France MSME summary indicators with more precise (2-decimal) figures.
# These cells hold their figures as text (not numbers), so the target number
# format is forced to Text before writing the value - otherwise Excel would
# helpfully "fix" the look-alike numeric string back into a real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Enterprises density (per 1000 people) - row 13 (Micro / SMEs / MSMEs)
Set-TextValue "B13" "36.98"
Set-TextValue "C13" "2.22"
Set-TextValue "D13" "39.21"

# Employment (% of total) - row 14 (Micro / SMEs / MSMEs)
Set-TextValue "B14" "28.65"
Set-TextValue "C14" "33.96"
Set-TextValue "D14" "62.61"

# Enterprises (% of total) - row 16 (Micro / SMEs / MSMEs)
Set-TextValue "B16" "94.17"
Set-TextValue "C16" "5.66"
Set-TextValue "D16" "99.83"
